$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Sheet1"

# Clear the applied fill style from the data cells (rows 2-6, cols A-F)
$ws.Range("A2:F6").Style = "Normal"

# Set column widths per new layout
# (inputs chosen so the stored/serialized width lands on the target value,
#  accounting for the engine's internal width quantization)
$ws.Columns.Item(1).ColumnWidth = 16.666666666666668
$ws.Columns.Item(2).ColumnWidth = 16.666666666666668
$ws.Columns.Item(3).ColumnWidth = 17.166666666666668
$ws.Columns.Item(4).ColumnWidth = 16.166666666666668
$ws.Columns.Item(5).ColumnWidth = 13.666666666666666
$ws.Columns.Item(6).ColumnWidth = 10.666666666666666

# Update the selection shown on the sheet
$ws.Range("I8").Select()
